{"js": "// Replace each two-digit multiplication problem's text with its new value.\n// Each \"AA\u00d7BB=\" string is unique in the document, so a scoped, case-sensitive\n// search-and-replace for each old/new pair safely targets the right run.\nconst replacements = [\n  [\"45\u00d754=\", \"13\u00d781=\"],\n  [\"64\u00d718=\", \"92\u00d738=\"],\n  [\"57\u00d798=\", \"62\u00d781=\"],\n  [\"33\u00d744=\", \"92\u00d783=\"],\n  [\"59\u00d757=\", \"42\u00d756=\"],\n  [\"40\u00d724=\", \"85\u00d756=\"],\n  [\"50\u00d737=\", \"44\u00d718=\"],\n  [\"88\u00d753=\", \"43\u00d771=\"],\n  [\"28\u00d792=\", \"32\u00d794=\"],\n  [\"60\u00d782=\", \"41\u00d763=\"],\n  [\"12\u00d771=\", \"17\u00d785=\"],\n  [\"95\u00d754=\", \"51\u00d787=\"],\n  [\"25\u00d769=\", \"86\u00d714=\"],\n  [\"75\u00d774=\", \"99\u00d778=\"],\n  [\"19\u00d718=\", \"29\u00d756=\"],\n  [\"73\u00d727=\", \"50\u00d786=\"],\n  [\"64\u00d719=\", \"44\u00d761=\"],\n  [\"80\u00d749=\", \"35\u00d744=\"],\n  [\"25\u00d777=\", \"65\u00d783=\"],\n  [\"63\u00d761=\", \"19\u00d786=\"],\n  [\"74\u00d784=\", \"19\u00d784=\"],\n  [\"68\u00d769=\", \"56\u00d763=\"],\n  [\"24\u00d785=\", \"26\u00d732=\"],\n  [\"74\u00d762=\", \"43\u00d785=\"],\n  [\"11\u00d716=\", \"17\u00d748=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication problem's text with its new value.\n# Each \"AA\u00d7BB=\" string is unique in the document, so a case-sensitive\n# whole-document Find/Replace for each old/new pair safely targets the\n# correct run without disturbing any other content or formatting.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"45\u00d754=\", \"13\u00d781=\"),\n    @(\"64\u00d718=\", \"92\u00d738=\"),\n    @(\"57\u00d798=\", \"62\u00d781=\"),\n    @(\"33\u00d744=\", \"92\u00d783=\"),\n    @(\"59\u00d757=\", \"42\u00d756=\"),\n    @(\"40\u00d724=\", \"85\u00d756=\"),\n    @(\"50\u00d737=\", \"44\u00d718=\"),\n    @(\"88\u00d753=\", \"43\u00d771=\"),\n    @(\"28\u00d792=\", \"32\u00d794=\"),\n    @(\"60\u00d782=\", \"41\u00d763=\"),\n    @(\"12\u00d771=\", \"17\u00d785=\"),\n    @(\"95\u00d754=\", \"51\u00d787=\"),\n    @(\"25\u00d769=\", \"86\u00d714=\"),\n    @(\"75\u00d774=\", \"99\u00d778=\"),\n    @(\"19\u00d718=\", \"29\u00d756=\"),\n    @(\"73\u00d727=\", \"50\u00d786=\"),\n    @(\"64\u00d719=\", \"44\u00d761=\"),\n    @(\"80\u00d749=\", \"35\u00d744=\"),\n    @(\"25\u00d777=\", \"65\u00d783=\"),\n    @(\"63\u00d761=\", \"19\u00d786=\"),\n    @(\"74\u00d784=\", \"19\u00d784=\"),\n    @(\"68\u00d769=\", \"56\u00d763=\"),\n    @(\"24\u00d785=\", \"26\u00d732=\"),\n    @(\"74\u00d762=\", \"43\u00d785=\"),\n    @(\"11\u00d716=\", \"17\u00d748=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
